# Apply row-permutation updates for Hortaliza / Rabanito dataset ("Fruta / hortaliza, semanal")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44161
$ws.Range("J2").Value = 7000
$ws.Range("D3").Value = 44214
$ws.Range("D4").Value = 44188
$ws.Range("D5").Value = 44600
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 1300
$ws.Range("K5").Value = 3500
$ws.Range("L5").Value = 4000
$ws.Range("M5").Value = 3808
$ws.Range("O5").Value = 'Región Metropolitana'
$ws.Range("P5").Value = 38
$ws.Range("D6").Value = 44186
$ws.Range("J6").Value = 10000
$ws.Range("D7").Value = 44187
$ws.Range("J7").Value = 12000
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 3000
$ws.Range("O7").Value = 'Provincia de Chacabuco'
$ws.Range("P7").Value = 30
$ws.Range("D8").Value = 44160
$ws.Range("J8").Value = 7000
$ws.Range("D9").Value = 44204
$ws.Range("J9").Value = 7000
$ws.Range("D11").Value = 44159
$ws.Range("J11").Value = 7000
$ws.Range("D12").Value = 44210
$ws.Range("J12").Value = 8800
$ws.Range("K12").Value = 2500
$ws.Range("M12").Value = 2750
$ws.Range("O12").Value = 'Provincia de Chacabuco'
$ws.Range("P12").Value = 28
$ws.Range("D13").Value = 44602
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 12000
$ws.Range("K13").Value = 3000
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = 3000
$ws.Range("O13").Value = 'Provincia de Chacabuco'
$ws.Range("P13").Value = 30
$ws.Range("D14").Value = 44602
$ws.Range("I14").Value = 'Segunda'
$ws.Range("J14").Value = 6000
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2500
$ws.Range("P14").Value = 25
$ws.Range("D15").Value = 44166
$ws.Range("D16").Value = 44181
$ws.Range("J16").Value = 12000
$ws.Range("D17").Value = 44229
$ws.Range("J17").Value = 16000
$ws.Range("D18").Value = 44167
$ws.Range("J18").Value = 7000
$ws.Range("K18").Value = 3000
$ws.Range("M18").Value = 3000
$ws.Range("P18").Value = 30
$ws.Range("D19").Value = 44230
$ws.Range("J19").Value = 16000
$ws.Range("D20").Value = 44245
$ws.Range("J20").Value = 9000
$ws.Range("O20").Value = 'Región Metropolitana'
$ws.Range("D21").Value = 44245
$ws.Range("I21").Value = 'Segunda'
$ws.Range("J21").Value = 5000
$ws.Range("K21").Value = 2500
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = 2500
$ws.Range("O21").Value = 'Región Metropolitana'
$ws.Range("P21").Value = 25
$ws.Range("D23").Value = 44162
$ws.Range("J23").Value = 7000
$ws.Range("D24").Value = 44168
$ws.Range("J24").Value = 7000
$ws.Range("D25").Value = 44232
$ws.Range("J25").Value = 16000
$ws.Range("D26").Value = 44215
$ws.Range("J26").Value = 16000
$ws.Range("D27").Value = 44189
$ws.Range("J27").Value = 16000
